$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1148.0526
$ws.Range("I28").Value = 1242.3
$ws.Range("J28").Value = 1043.3334
$ws.Range("K28").Value = 1242.3
$ws.Range("L28").Value = 1043.3334
$ws.Range("M28").Value = -757.3
$ws.Range("N28").Value = -2013.3334
$ws.Range("H69").Value = 5500
$ws.Range("J69").Value = 5500
$ws.Range("L69").Value = 16500
$ws.Range("N69").Value = -18248
$ws.Range("H72").Value = 5500
$ws.Range("J72").Value = 5500
$ws.Range("L72").Value = 49500
$ws.Range("N72").Value = -58236
$ws.Range("H76").Value = 3148.2415
$ws.Range("I76").Value = 3153.5356
$ws.Range("J76").Value = 3000
$ws.Range("K76").Value = 3153.5356
$ws.Range("L76").Value = 3000
$ws.Range("M76").Value = -2838.5356
$ws.Range("N76").Value = -3630
$ws.Range("H79").Value = 3148.2415
$ws.Range("I79").Value = 3153.5356
$ws.Range("J79").Value = 3000
$ws.Range("K79").Value = 3153.5356
$ws.Range("L79").Value = 3000
$ws.Range("M79").Value = -2061.5356
$ws.Range("N79").Value = -5184
$ws.Range("H107").Value = 1311.3684
$ws.Range("I107").Value = 1483.8572
$ws.Range("J107").Value = 828.4
$ws.Range("K107").Value = 1483.8572
$ws.Range("L107").Value = 828.4
$ws.Range("M107").Value = 436.1428000000001
$ws.Range("N107").Value = -4668.4
$ws.Range("H112").Value = 23810722
$ws.Range("J112").Value = 24391458
$ws.Range("L112").Value = 73174374
$ws.Range("N112").Value = -73176590
$ws.Range("H115").Value = 1052.8572
$ws.Range("I115").Value = 342.5
$ws.Range("J115").Value = 2000
$ws.Range("K115").Value = 1027.5
$ws.Range("L115").Value = 6000
$ws.Range("M115").Value = 539.5
$ws.Range("N115").Value = -9134
$ws.Range("H132").Value = 288633.2
$ws.Range("I132").Value = 373791
$ws.Range("K132").Value = 1121373
$ws.Range("M132").Value = -1118843
$ws.Range("H137").Value = 22223644
$ws.Range("I137").Value = 1000.0769
$ws.Range("J137").Value = 52633580
$ws.Range("K137").Value = 3000.2307
$ws.Range("L137").Value = 157900740
$ws.Range("M137").Value = -450.2307000000001
$ws.Range("N137").Value = -157905840

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 2749.9614
$ws.Range("H132").Value = 1577.5957
$ws.Range("I132").Value = 1461.9642
$ws.Range("J132").Value = 1748
$ws.Range("K132").Value = 4385.892599999999
$ws.Range("L132").Value = 5244
$ws.Range("M132").Value = -1855.892599999999
$ws.Range("N132").Value = -10304

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1010.2759
$ws.Range("I20").Value = 730.5
$ws.Range("J20").Value = 1632
$ws.Range("K20").Value = 730.5
$ws.Range("L20").Value = 1632
$ws.Range("M20").Value = -483.5
$ws.Range("N20").Value = -2126
$ws.Range("H94").Value = 1268.0555
$ws.Range("I94").Value = 1092.1666
$ws.Range("J94").Value = 1619.8334
$ws.Range("K94").Value = 1092.1666
$ws.Range("L94").Value = 1619.8334
$ws.Range("M94").Value = -641.1666
$ws.Range("N94").Value = -2521.8334

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 8000
$ws.Range("J4").Value = 8000
$ws.Range("L4").Value = 8000
$ws.Range("N4").Value = -8224
$ws.Range("H31").Value = 2333.6956
$ws.Range("I31").Value = 1106.8182
$ws.Range("J31").Value = 3458.3333
$ws.Range("K31").Value = 1106.8182
$ws.Range("L31").Value = 3458.3333
$ws.Range("M31").Value = -811.8181999999999
$ws.Range("N31").Value = -4048.3333
$ws.Range("H34").Value = 2333.6956
$ws.Range("I34").Value = 1106.8182
$ws.Range("J34").Value = 3458.3333
$ws.Range("K34").Value = 1106.8182
$ws.Range("L34").Value = 3458.3333
$ws.Range("M34").Value = -904.8181999999999
$ws.Range("N34").Value = -3862.3333

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 876580.7
$ws.Range("I64").Value = 968.5
$ws.Range("J64").Value = 1401948
$ws.Range("K64").Value = 2905.5
$ws.Range("L64").Value = 4205844
$ws.Range("M64").Value = -2635.5
$ws.Range("N64").Value = -4206384
$ws.Range("H67").Value = 876580.7
$ws.Range("I67").Value = 968.5
$ws.Range("J67").Value = 1401948
$ws.Range("K67").Value = 2905.5
$ws.Range("L67").Value = 4205844
$ws.Range("M67").Value = -1969.5
$ws.Range("N67").Value = -4207716
$ws.Range("H117").Value = 1826.6666
$ws.Range("I117").Value = 1250
$ws.Range("K117").Value = 3750
$ws.Range("M117").Value = -308
$ws.Range("H122").Value = 22224040
$ws.Range("I122").Value = 55556132
$ws.Range("J122").Value = 2646
$ws.Range("K122").Value = 500005188
$ws.Range("L122").Value = 23814
$ws.Range("M122").Value = -500002738
$ws.Range("N122").Value = -28714
$ws.Range("H129").Value = 928.1818
$ws.Range("I129").Value = 527.8
$ws.Range("J129").Value = 1261.8334
$ws.Range("K129").Value = 1583.4
$ws.Range("L129").Value = 3785.5002
$ws.Range("M129").Value = 3416.6
$ws.Range("N129").Value = -13785.5002
$ws.Range("H132").Value = 50002550
$ws.Range("I132").Value = 125002090
$ws.Range("J132").Value = 2866.6667
$ws.Range("K132").Value = 1125018810
$ws.Range("L132").Value = 25800.0003
$ws.Range("M132").Value = -1125016280
$ws.Range("N132").Value = -30860.0003

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 2319.7222
$ws.Range("J5").Value = 2319.7222
$ws.Range("L5").Value = 2319.7222
$ws.Range("N5").Value = -2543.7222
$ws.Range("H70").Value = 4192.409
$ws.Range("I70").Value = 3864.923
$ws.Range("K70").Value = 3864.923
$ws.Range("M70").Value = -3594.923
$ws.Range("H73").Value = 4192.409
$ws.Range("I73").Value = 3864.923
$ws.Range("K73").Value = 3864.923
$ws.Range("M73").Value = -2928.923
$ws.Range("H80").Value = 121981.1
$ws.Range("I80").Value = 2402.5
$ws.Range("J80").Value = 151875.75
$ws.Range("K80").Value = 2402.5
$ws.Range("L80").Value = 151875.75
$ws.Range("M80").Value = -1404.5
$ws.Range("N80").Value = -153871.75
$ws.Range("H83").Value = 121981.1
$ws.Range("I83").Value = 2402.5
$ws.Range("J83").Value = 151875.75
$ws.Range("K83").Value = 12012.5
$ws.Range("L83").Value = 759378.75
$ws.Range("M83").Value = -7020.5
$ws.Range("N83").Value = -769362.75
$ws.Range("H93").Value = 20000
$ws.Range("J93").Value = 20000
$ws.Range("L93").Value = 20000
$ws.Range("N93").Value = -23744
$ws.Range("H138").Value = 39000
$ws.Range("J138").Value = 39000
$ws.Range("L138").Value = 39000
$ws.Range("N138").Value = -49280

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 10336001
$ws.Range("J2").Value = 15004001
$ws.Range("L2").Value = 15004001
$ws.Range("N2").Value = -15004225
$ws.Range("H62").Value = 13250
$ws.Range("I62").Value = 9800
$ws.Range("J62").Value = 16700
$ws.Range("K62").Value = 9800
$ws.Range("L62").Value = 16700
$ws.Range("M62").Value = -9176
$ws.Range("N62").Value = -17948
$ws.Range("H64").Value = 11100
$ws.Range("J64").Value = 11100
$ws.Range("L64").Value = 11100
$ws.Range("N64").Value = -11550
$ws.Range("H65").Value = 13250
$ws.Range("I65").Value = 9800
$ws.Range("J65").Value = 16700
$ws.Range("K65").Value = 29400
$ws.Range("L65").Value = 50100
$ws.Range("M65").Value = -26280
$ws.Range("N65").Value = -56340
$ws.Range("H67").Value = 11100
$ws.Range("J67").Value = 11100
$ws.Range("L67").Value = 11100
$ws.Range("N67").Value = -12660
$ws.Range("H106").Value = 11033.167
$ws.Range("J106").Value = 11033.167
$ws.Range("L106").Value = 11033.167
$ws.Range("N106").Value = -13557.167

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 40000000
$ws.Range("I2").Value = 10000000
$ws.Range("K2").Value = 10000000
$ws.Range("H96").Value = 1281.5714
$ws.Range("J96").Value = 1394.2
$ws.Range("L96").Value = 1394.2
$ws.Range("N96").Value = -4140.2
$ws.Range("H104").Value = 33950
$ws.Range("J104").Value = 33950
$ws.Range("L104").Value = 33950
$ws.Range("N104").Value = -40938
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 449.91666
$ws.Range("I107").Value = 267.7143
$ws.Range("J107").Value = 705
$ws.Range("K107").Value = 803.1428999999999
$ws.Range("L107").Value = 2115
$ws.Range("M107").Value = 1116.8571
$ws.Range("N107").Value = -5955
